# Apply the daily crypto price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.092.71'
$ws.Range("E2").Value = '  -2.87%  '

# Row 3
$ws.Range("D3").Value = '1.843.43'

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = "'0.6983"
$ws.Range("E5").Value = '  -5.79%  '

# Row 6
$ws.Range("D6").Value = "'237.36"
$ws.Range("E6").Value = '  -2.27%  '

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = "'0.3026"
$ws.Range("E8").Value = '  -4.26%  '

# Row 9
$ws.Range("D9").Value = "'0.07412"
$ws.Range("E9").Value = '  +2.47%  '

# Row 10
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = '  -6.87%  '

# Row 11
$ws.Range("D11").Value = "'0.08103"
$ws.Range("E11").Value = '  -2.98%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.841.61'
$ws.Range("E12").Value = '  -4.69%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = "'0.7225"
$ws.Range("E13").Value = '  -4.61%  '

# Row 14
$ws.Range("D14").Value = "'5.204"
$ws.Range("E14").Value = '  -3.80%  '

# Row 15
$ws.Range("D15").Value = "'88.84"
$ws.Range("E15").Value = '  -4.04%  '

# Row 16
$ws.Range("D16").Value = '29.018.50'
$ws.Range("E16").Value = '  -3.14%  '

# Row 17
$ws.Range("D17").Value = "'5.793"
$ws.Range("E17").Value = '  -5.83%  '

# Row 18
$ws.Range("D18").Value = "'240.87"
$ws.Range("E18").Value = '  -3.48%  '

# Row 19
$ws.Range("D19").Value = "'0.000007643"
$ws.Range("E19").Value = '  -2.78%  '

# Row 20
$ws.Range("E20").Value = '  -4.63%  '

# Row 21
$ws.Range("D21").Value = "'0.9988"
$ws.Range("E21").Value = '  -0.14%  '

# Row 22
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  +0.01%  '

# Row 23
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.073.42'
$ws.Range("E23").Value = '  -3.81%  '

# Row 24
$ws.Range("D24").Value = "'7.544"
$ws.Range("E24").Value = '  -5.94%  '

# Row 25
$ws.Range("D25").Value = "'0.1475"
$ws.Range("E25").Value = '  -5.81%  '

# Row 26
$ws.Range("D26").Value = "'162.06"
$ws.Range("E26").Value = '  -2.10%  '

# Row 27
$ws.Range("D27").Value = "'8.939"
$ws.Range("E27").Value = '  -4.09%  '

# Row 28
$ws.Range("D28").Value = "'17.98"
$ws.Range("E28").Value = '  -3.98%  '

# Row 29
$ws.Range("D29").Value = "'1.933"
$ws.Range("E29").Value = '  -5.15%  '

# Row 30
$ws.Range("D30").Value = "'1.373"
$ws.Range("E30").Value = '  -8.02%  '

# Row 31
$ws.Range("D31").Value = "'4.449"
$ws.Range("E31").Value = '  -3.72%  '

# Row 32
$ws.Range("E32").Value = '  -3.14%  '

# Row 33
$ws.Range("D33").Value = "'4.010"
$ws.Range("E33").Value = '  -5.30%  '

# Row 34
$ws.Range("D34").Value = "'0.05183"
$ws.Range("E34").Value = '  -3.69%  '

# Row 35
$ws.Range("D35").Value = "'1.182"
$ws.Range("E35").Value = '  -5.64%  '

# Row 36
$ws.Range("D36").Value = "'0.7095"
$ws.Range("E36").Value = '  -6.39%  '

# Row 37
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  -0.31%  '

# Row 38
$ws.Range("D38").Value = "'2.647"
$ws.Range("E38").Value = '  -2.33%  '

# Row 39
$ws.Range("D39").Value = "'0.01870"
$ws.Range("E39").Value = '  -5.25%  '

# Row 40
$ws.Range("D40").Value = "'2.673"
$ws.Range("E40").Value = '  -3.20%  '

# Row 41
$ws.Range("D41").Value = "'0.9022"
$ws.Range("E41").Value = '  +4.64%  '

# Row 42
$ws.Range("D42").Value = "'0.4283"
$ws.Range("E42").Value = '  -6.22%  '

# Row 43
$ws.Range("D43").Value = "'5.918"
$ws.Range("E43").Value = '  -3.19%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.048.68'
$ws.Range("E44").Value = '  -5.10%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'69.86"
$ws.Range("E45").Value = '  -3.93%  '

# Row 46
$ws.Range("D46").Value = "'0.9996"
$ws.Range("E46").Value = '  -0.12%  '

# Row 47
$ws.Range("D47").Value = "'101.54"
$ws.Range("E47").Value = '  -3.14%  '

# Row 48
$ws.Range("D48").Value = "'1.749"
$ws.Range("E48").Value = '  -6.56%  '

# Row 49
$ws.Range("D49").Value = "'7.078"
$ws.Range("E49").Value = '  -7.19%  '

# Row 50
$ws.Range("D50").Value = "'9.179"
$ws.Range("E50").Value = '  -3.98%  '

# Row 51
$ws.Range("D51").Value = '1.978.07'
$ws.Range("E51").Value = '  -4.76%  '
